$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column for rows 2-6 from 2023-09-06 to 2023-09-14
$newDate = Get-Date -Year 2023 -Month 9 -Day 14 -Hour 0 -Minute 0 -Second 0
$newDate = $newDate.Date
$ws.Range("C2:C6").Value = $newDate
